# Update test data for removed variables and rely more on moorings &
# installation examples.
#
# The "On-Site" sheet's spare-parts table only listed four sub-systems
# (Inter-Array Cables, Substations, Export Cable + header). This adds a
# fifth row of data for the "Umbilical Cable" sub-system, reusing the
# same formatting as the rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("On-Site")

# Fill in the previously-empty row 5 with the new sub-system's data.
$ws.Range("A5").Value = "Umbilical Cable"
$ws.Range("B5").Value = 30000
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 6

# Row 5 already carried the "Normal 2 + fill/border" style used for the
# rest of column A; bring the other data rows' A-column cells onto the
# same style so the whole column is consistent.
$ws.Range("A2:A5").Style = $ws.Range("A5").Style

# Leave the selection on the newly-completed cell, as the author did.
[void]$ws.Range("A5").Select()
